$wb = $excel.ActiveWorkbook

# --- Rename existing sheet, add new sheet "Distribution" right after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Goal-oriented trivial"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Distribution"

# --- Update a couple of values on sheet1 ---
$ws1.Range("C6").Value = 0.94589
$ws1.Range("C7").Value = 0.9896

# --- Fix selection on sheet1 (no longer the active/tab-selected sheet) ---
$ws1.Range("C14").Select()

# --- Populate sheet2 "Distribution" ---
# Row labels first (column A), matching the original authoring order so
# shared-string indices line up with the target file.
$ws2.Range("A2").Value = "Validity"
$ws2.Range("A3").Value = "Uniqueness"
$ws2.Range("A4").Value = "Novelty"
$ws2.Range("A5").Value = "KL divergence"
$ws2.Range("A6").Value = "Frechet ChemNet distance"

# Column headers (row 1), written E1, B1, C1, D1 to match original order
$ws2.Range("E1").Value = "With discriminator"
$ws2.Range("B1").Value = "No priors"
$ws2.Range("C1").Value = "Unconditional priors"
$ws2.Range("D1").Value = "Conditional priors"

# Numeric data
$ws2.Range("D2").Value = 1
$ws2.Range("E2").Value = 1

$ws2.Range("D3").Value = 1
$ws2.Range("E3").Value = 0.9987

$ws2.Range("D4").Value = 0.9873
$ws2.Range("E4").Value = 0.9946

$ws2.Range("D5").Value = 0.73
$ws2.Range("E5").Value = 0.8

$ws2.Range("D6").Value = 0.0531
$ws2.Range("E6").Value = 0.11

# Column widths (closest achievable inputs given the host's width-rounding
# model; targets are 22.47265625 / 8.89453125 / 9.3125 character-widths)
$ws2.Range("A1").ColumnWidth = 21.666666666666668
$ws2.Range("B1").ColumnWidth = 8
$ws2.Range("C1").ColumnWidth = 8.5

# --- Selection + active sheet on the Distribution sheet ---
$ws2.Range("D2").Select()
$ws2.Activate()
